# Add the "invalid number of offers" row (row 4) to the sheet, mirroring
# rows 1-3's layout but extended out to column I.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = "ABCHJUH"
$ws.Range("B4").Value = 8
$ws.Range("C4").Value = "909ikokujyhtgt*"
$ws.Range("D4").Value = "JKJKUHY/////\\\\\%^%gyvb"
$ws.Range("E4").Value = "iojkjkjhjhjhjjhj"
$ws.Range("F4").Value = "ioiojkhjghfgfgghbn"
$ws.Range("G4").Value = "jkjkhjhh)))))"
$ws.Range("H4").Value = "(((hjnmnmnmm####"
$ws.Range("I4").Value = "hjhj.uiuiuisdksd"

# Best-fit the newly used column H, same as the other data columns.
$ws.Columns.Item(8).ColumnWidth = 16.721354166666668

# Selection moves to the new empty column G (whole column) as left in the
# saved workbook.
$ws.Range("G1:G1048576").Select() | Out-Null
